# Update cryptos list (price + volume%) and fix the ONDO/Stacks row order,
# as scraped by the GitHub Actions job.
# NumberFormat "@" is applied before writing any Price cell whose new text
# looks like a plain number (e.g. "1.00", "7.03") so Excel stores it as text
# instead of silently coercing it to a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.518.60'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '3.549.53'
$ws.Range("E3").Value = '  +3.14%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.06'
$ws.Range("E5").Value = '  +2.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.85'
$ws.Range("E6").Value = '  +3.08%  '
$ws.Range("D7").Value = '3.550.89'
$ws.Range("E7").Value = '  +3.20%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("E10").Value = '  +2.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.03'
$ws.Range("E11").Value = '  -4.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.395'
$ws.Range("E12").Value = '  +4.41%  '
$ws.Range("D13").Value = '4.154.47'
$ws.Range("E13").Value = '  +3.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000188'
$ws.Range("E14").Value = '  +3.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.31'
$ws.Range("E15").Value = '  +2.44%  '
$ws.Range("D16").Value = '3.561.09'
$ws.Range("E16").Value = '  +3.30%  '
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("D18").Value = '65.560.17'
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.31'
$ws.Range("E19").Value = '  +3.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.94'
$ws.Range("E20").Value = '  +2.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.34'
$ws.Range("E21").Value = '  +4.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '396.33'
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.573'
$ws.Range("E23").Value = '  +4.55%  '
$ws.Range("D24").Value = '3.698.12'
$ws.Range("E24").Value = '  +3.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.29'
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  +11.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.94'
$ws.Range("E28").Value = '  +9.88%  '
$ws.Range("E29").Value = '  +2.72%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.37'
$ws.Range("E31").Value = '  +2.02%  '
$ws.Range("D32").Value = '3.565.05'
$ws.Range("E32").Value = '  +3.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.148'
$ws.Range("E33").Value = '  +1.18%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.82'
$ws.Range("E35").Value = '  +3.52%  '
$ws.Range("E36").Value = '  +8.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.06'
$ws.Range("E37").Value = '  +1.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '169.54'
$ws.Range("E38").Value = '  -1.14%  '
$ws.Range("E39").Value = '  +3.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.07'
$ws.Range("E40").Value = '  +4.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0818'
$ws.Range("E41").Value = '  +6.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.835'
$ws.Range("E42").Value = '  +1.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.47'
$ws.Range("E43").Value = '  +16.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.99'
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.45'
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.71'
$ws.Range("E47").Value = '  +5.12%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.21'
$ws.Range("E48").Value = '  +9.79%  '
$ws.Range("D49").Value = '2.481.02'
$ws.Range("E49").Value = '  +12.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.83'
$ws.Range("E50").Value = '  +4.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.36'
$ws.Range("E51").Value = '  +19.11%  '
